$d = $word.ActiveDocument

# Replace the original text "Prueba 1." with the new wording.
$d.Content.Find.Execute("Prueba 1.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Esta es la versión modificada por Guillermo Ruiz Vida", 2)

# Append 15 blank paragraphs after the (now updated) first paragraph, right
# before the section properties. Collapsing the range to its end and typing
# a single carriage return at a time (rather than several at once) produces
# clean, content-less <w:p/> paragraphs instead of ones carrying an empty run.
$r = $d.Content
for ($i = 0; $i -lt 15; $i++) {
    $r.Collapse(0)
    $r.Text = "`r"
}
